# Update the build timestamp embedded in the version strings from
# "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet ---
$aboutSheet = $wb.Worksheets.Item("About")

$aboutSheet.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$aboutSheet.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Shanxi Jinyuan Coal Mine, China, M0330, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)" + [char]34

# --- "Boundaries and methane sources" sheet ---
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 14; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19) # column S
    $cell.Value = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
}
